$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Daniel Sams"

# Insert a new column before column A (shifts teamName..result from A:L to B:M)
$ws.Columns.Item(1).Insert()

# Insert a new row before the current data row 2 (pushes the Punjab Kings match to row 3)
$ws.Rows.Item(2).Insert()

# Keep everything as text, matching the original sheet's stored-as-text values
$ws.Range("A1:M3").NumberFormat = "@"

# Header for the new column
$ws.Range("A1").Value = "matchNo"

# matchNo for the existing (now row 3) match
$ws.Range("A3").Value = "26th"

# New row 2: the 22nd match vs Delhi Capitals
$ws.Range("A2").Value = "22nd"
$ws.Range("B2").Value = "Royal Challengers Bangalore"
$ws.Range("C2").Value = "Daniel Sams"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "3"
$ws.Range("F2").Value = "2"
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "150.00"
$ws.Range("J2").Value = "Delhi Capitals"
$ws.Range("K2").Value = "Ahmedabad"
$ws.Range("L2").Value = "April 27"
$ws.Range("M2").Value = "RCB won by 1 run"
